$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PatternIndex (column A) text for rows 10-19.
# This mirrors the author's edit to the trigram model: the last five
# pattern names were replaced with a new batch of organizational
# patterns, and the remaining names in A10:A14 shifted up to absorb the
# removed "OrganizationFollowsMarket" entry.
$ws.Range("A10").Value = "DevelopingInPairs"
$ws.Range("A11").Value = "DistributeWorkEvenly"
$ws.Range("A12").Value = "DivideAndConquer"
$ws.Range("A13").Value = "DomainExpertiseInRoles"
$ws.Range("A14").Value = "FeatureAssignment"
$ws.Range("A15").Value = "Stand-UpMeeting"
$ws.Range("A16").Value = "StandardsLinkingLocations"
$ws.Range("A17").Value = "SubclassPerTeam"
$ws.Range("A18").Value = "TeamPerTask"
$ws.Range("A19").Value = "VariationBehindInterface"

# Move the active selection to the freshly edited block.
$ws.Range("A10:M19").Select()
